$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193 (pushes existing rows 193:235 down to 194:236)
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new weekly "Ajo" price record
$ws.Range("A193").Value = 4
$ws.Range("B193").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C193").Value = "Los Lagos"
$ws.Range("D193").Value = 44641
$ws.Range("E193").Value = 10
$ws.Range("F193").Value = 100112003
$ws.Range("G193").Value = "Ajo"
$ws.Range("H193").Value = "Chino"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 80
$ws.Range("K193").Value = 21000
$ws.Range("L193").Value = 21000
$ws.Range("M193").Value = 21000
$ws.Range("N193").Value = "`$/caja 10 kilos"
$ws.Range("O193").Value = "China"
$ws.Range("P193").Value = 2100
$ws.Range("Q193").Value = 10
$ws.Range("R193").Value = "Hortaliza"
